# The post "「ゆっくりと」بشويش" (row 815) was removed from the sheet.
# Deleting the entire row shifts every subsequent row up by one, which
# matches the diff (rows 816-837 become 815-836) and shrinks the used
# range from A1:C837 to A1:C836.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(815).Delete()
